$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19 (shifts CORP_CARD and everything below it down by one)
$ws.Rows("19:19").Insert()

# The inserted row doesn't automatically inherit the bordered-table formatting,
# so copy it over explicitly: column A formatting from the row now below (old row 19),
# and column B's "highlighted" style from another B cell that already used it.
$ws.Range("A20:B20").Copy() | Out-Null
$ws.Range("A19:B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null

# Fill in the new indicator row
$ws.Range("A19").Value = "SP_BIG_PLANS"
$ws.Range("B19").Value = "Доля клиентов с ПУ Большие планы"

# Match the author's final selection/scroll position
$ws.Range("B19").Select() | Out-Null
